# "random updates to sync"
# Duplicate the "Leafouts" sheet, place the copy before the original, name it
# "Leafouts Rand", and leave the cursor/selection on each sheet the way the
# author left them: K16 on the new random-sample copy, C19 on the original
# "Leafouts" sheet (which stays the active tab). "Nonleafouts" is untouched.

$wb = $excel.ActiveWorkbook

$leafouts = $wb.Worksheets.Item("Leafouts")

# Create a copy of "Leafouts" and drop it immediately before "Leafouts" itself.
$leafouts.Copy($leafouts)

# The freshly created copy is now the first sheet in the workbook.
$leafoutsRand = $wb.Worksheets.Item(1)
$leafoutsRand.Name = "Leafouts Rand"

# Set the lingering selection on the new sheet first (without activating it),
# then activate "Leafouts" and set its selection last so it ends up as the
# active/selected tab. Re-fetch "Leafouts" by name: after Copy() shifted sheet
# positions, the old $leafouts handle would otherwise still resolve to the
# sheet that is now at index 1 (the copy), not the original sheet.
$leafoutsRand.Range("K16").Select()

$leafouts = $wb.Worksheets.Item("Leafouts")
$leafouts.Activate()
$leafouts.Range("C19").Select()
